{"js": "const replacements = [\n  [\"228\u00d76=\", \"719\u00d72=\"],\n  [\"924\u00d76=\", \"351\u00d77=\"],\n  [\"732\u00d74=\", \"395\u00d78=\"],\n  [\"146\u00d77=\", \"220\u00d74=\"],\n  [\"370\u00d72=\", \"387\u00d76=\"],\n  [\"374\u00d73=\", \"999\u00d78=\"],\n  [\"476\u00d75=\", \"617\u00d78=\"],\n  [\"623\u00d73=\", \"830\u00d73=\"],\n  [\"396\u00d73=\", \"660\u00d75=\"],\n  [\"510\u00d72=\", \"662\u00d77=\"],\n  [\"221\u00d75=\", \"141\u00d72=\"],\n  [\"432\u00d72=\", \"330\u00d73=\"],\n  [\"763\u00d77=\", \"969\u00d73=\"],\n  [\"559\u00d78=\", \"985\u00d73=\"],\n  [\"267\u00d75=\", \"223\u00d74=\"],\n  [\"285\u00d79=\", \"361\u00d78=\"],\n  [\"226\u00d76=\", \"604\u00d78=\"],\n  [\"908\u00d78=\", \"216\u00d76=\"],\n  [\"710\u00d75=\", \"794\u00d78=\"],\n  [\"474\u00d75=\", \"268\u00d76=\"],\n  [\"883\u00d73=\", \"168\u00d72=\"],\n  [\"175\u00d74=\", \"259\u00d73=\"],\n  [\"740\u00d72=\", \"353\u00d78=\"],\n  [\"303\u00d72=\", \"393\u00d77=\"],\n  [\"376\u00d79=\", \"592\u00d74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"228\u00d76=\", \"719\u00d72=\")\n  ,@(\"924\u00d76=\", \"351\u00d77=\")\n  ,@(\"732\u00d74=\", \"395\u00d78=\")\n  ,@(\"146\u00d77=\", \"220\u00d74=\")\n  ,@(\"370\u00d72=\", \"387\u00d76=\")\n  ,@(\"374\u00d73=\", \"999\u00d78=\")\n  ,@(\"476\u00d75=\", \"617\u00d78=\")\n  ,@(\"623\u00d73=\", \"830\u00d73=\")\n  ,@(\"396\u00d73=\", \"660\u00d75=\")\n  ,@(\"510\u00d72=\", \"662\u00d77=\")\n  ,@(\"221\u00d75=\", \"141\u00d72=\")\n  ,@(\"432\u00d72=\", \"330\u00d73=\")\n  ,@(\"763\u00d77=\", \"969\u00d73=\")\n  ,@(\"559\u00d78=\", \"985\u00d73=\")\n  ,@(\"267\u00d75=\", \"223\u00d74=\")\n  ,@(\"285\u00d79=\", \"361\u00d78=\")\n  ,@(\"226\u00d76=\", \"604\u00d78=\")\n  ,@(\"908\u00d78=\", \"216\u00d76=\")\n  ,@(\"710\u00d75=\", \"794\u00d78=\")\n  ,@(\"474\u00d75=\", \"268\u00d76=\")\n  ,@(\"883\u00d73=\", \"168\u00d72=\")\n  ,@(\"175\u00d74=\", \"259\u00d73=\")\n  ,@(\"740\u00d72=\", \"353\u00d78=\")\n  ,@(\"303\u00d72=\", \"393\u00d77=\")\n  ,@(\"376\u00d79=\", \"592\u00d74=\")\n)\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.Text = $pair[1]\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
